$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.258919239044189
$ws.Range("B1").Value = 1.792936205863953
$ws.Range("C1").Value = 4.126768112182617
$ws.Range("D1").Value = 3.329569101333618
$ws.Range("E1").Value = 1.156816244125366
